$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input data corrections. Columns B/H/J/K hold shared formulas (running
#     totals / simple sums) that recompute automatically once the raw inputs
#     below change, so only the raw input cells need to be touched. ---

# Row 66: new positive cases count corrected from 4 to 5
$ws.Range("C66").Value = 5

# Row 441: new positive cases count corrected from 63 to 62
$ws.Range("C441").Value = 62

# Row 444: new positive cases count corrected from 38 to 37
$ws.Range("C444").Value = 37

# Row 446: new positive cases count corrected from 11 to 35
$ws.Range("C446").Value = 35

# Row 447: new positive cases count corrected from 1 to 34
$ws.Range("C447").Value = 34

# Row 448 (2021-05-18): this day's data had not been entered yet; fill it in.
$ws.Range("C448").Value = 4
$ws.Range("E448").Value = 6
$ws.Range("F448").Value = 5
$ws.Range("G448").Value = 19

# Columns L and M are formatted as Text ("@"). Writing a number straight into
# a Text-formatted cell would store it as text, which doesn't match the
# source data (plain numeric 0s there, same as every other row) - so flip the
# format to General just long enough to write the numeric literal, then
# restore the original Text format.
$ws.Range("L448").NumberFormat = "General"
$ws.Range("L448").Value = 0
$ws.Range("L448").NumberFormat = "@"

$ws.Range("M448").NumberFormat = "General"
$ws.Range("M448").Value = 0
$ws.Range("M448").NumberFormat = "@"
